$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header cell (H1) into the new header
# cells I1:J1 so they pick up the same bold/centered/bordered style (s=1).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Set new header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New I0 / IF data values for rows 2-71.
$ijData = @{
    2 = @(5,5)
    3 = @(8,8)
    4 = @(1,1)
    5 = @(9,9)
    6 = @(6,6)
    7 = @(7,7)
    8 = @(7,7)
    9 = @(9,9)
    10 = @(6,7)
    11 = @(6,6)
    12 = @(9,9)
    13 = @(8,8)
    14 = @(10,10)
    15 = @(7,7)
    16 = @(12,12)
    17 = @(7,7)
    18 = @(6,6)
    19 = @(8,8)
    20 = @(6,6)
    21 = @(7,7)
    22 = @(8,8)
    23 = @(6,6)
    24 = @(7,7)
    25 = @(8,8)
    26 = @(7,7)
    27 = @(7,7)
    28 = @(7,7)
    29 = @(8,8)
    30 = @(12,12)
    31 = @(7,7)
    32 = @(8,8)
    33 = @(8,8)
    34 = @(7,8)
    35 = @(4,5)
    36 = @(5,5)
    37 = @(10,10)
    38 = @(9,9)
    39 = @(8,8)
    40 = @(8,8)
    41 = @(7,7)
    42 = @(7,8)
    43 = @(8,8)
    44 = @(3,4)
    45 = @(6,7)
    46 = @(8,8)
    47 = @(7,8)
    48 = @(9,9)
    49 = @(6,8)
    50 = @(7,7)
    51 = @(7,8)
    52 = @(9,9)
    53 = @(6,6)
    54 = @(7,7)
    55 = @(7,7)
    56 = @(7,7)
    57 = @(7,7)
    58 = @(7,8)
    59 = @(7,7)
    60 = @(11,11)
    61 = @(5,6)
    62 = @(7,7)
    63 = @(9,9)
    64 = @(3,4)
    65 = @(8,8)
    66 = @(7,7)
    67 = @(8,9)
    68 = @(6,7)
    69 = @(9,9)
    70 = @(6,7)
    71 = @(3,3)
}

foreach ($r in $ijData.Keys) {
    $vals = $ijData[$r]
    $ws.Cells.Item($r, 9).Value2 = $vals[0]
    $ws.Cells.Item($r, 10).Value2 = $vals[1]
}
